# Revert "Fixed issues with form doubling up"
# Refreshes the Terminal Status / Cassette Balances report rows (5-27) with
# the prior snapshot's terminal IDs, names, balances, timestamps and cassette
# totals, and restores the prior grand total in E28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 580
$ws.Range("H5").Value = 45130.04188321759
$ws.Range("J5").Value = '07/17/23 18:04'
$ws.Range("K5").Value = '07/17/23 18:04'
$ws.Range("M5").Value = '$580 as of 7/17/2023 4:04:56 PM'
$ws.Range("N5").Value = 600
$ws.Range("A6").Value = 'L474746'
$ws.Range("C6").Value = 'ZACATES MARKET'
$ws.Range("E6").Value = 640
$ws.Range("H6").Value = 45129.04188321759
$ws.Range("J6").Value = '07/18/23 14:22'
$ws.Range("K6").Value = '07/18/23 14:22'
$ws.Range("M6").Value = '$780 as of 7/18/2023 10:19:11 AM'
$ws.Range("N6").Value = 680
$ws.Range("A7").Value = 'LK644532'
$ws.Range("C7").Value = 'SCL ENTERPRISES LAUNDRY'
$ws.Range("E7").Value = 700
$ws.Range("H7").Value = 45133.04188321759
$ws.Range("I7").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J7").Value = '07/17/23 20:34'
$ws.Range("K7").Value = '07/17/23 20:34'
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = '$700 as of 7/17/2023 6:34:04 PM'
$ws.Range("N7").Value = 760
$ws.Range("A9").Value = 'L704741'
$ws.Range("C9").Value = 'W ADAMS COIN LAUNDRY'
$ws.Range("E9").Value = 2100
$ws.Range("H9").Value = 45129.04188321759
$ws.Range("J9").Value = '07/18/23 15:07'
$ws.Range("K9").Value = '07/18/23 03:46'
$ws.Range("M9").Value = '$2,100 as of 7/18/2023 10:08:00 AM'
$ws.Range("N9").Value = 2100
$ws.Range("A10").Value = 'L662336'
$ws.Range("C10").Value = 'SB#4 MONA MARKET'
$ws.Range("E10").Value = 2260
$ws.Range("H10").Value = 45132.04188321759
$ws.Range("J10").Value = '07/17/23 17:03'
$ws.Range("K10").Value = '07/17/23 17:03'
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = '$2,260 as of 7/17/2023 3:03:55 PM'
$ws.Range("N10").Value = 2280
$ws.Range("A11").Value = 'LK864765'
$ws.Range("C11").Value = 'SKY LIQUOR'
$ws.Range("E11").Value = 2340
$ws.Range("H11").Value = 45130.04188321759
$ws.Range("J11").Value = '07/18/23 14:03'
$ws.Range("K11").Value = '07/18/23 10:02'
$ws.Range("M11").Value = '$2,420 as of 7/18/2023 8:02:45 AM'
$ws.Range("N11").Value = 2420
$ws.Range("A12").Value = 'LK561655'
$ws.Range("C12").Value = 'CRENSHAW CRAVOR #2'
$ws.Range("E12").Value = 2780
$ws.Range("H12").ClearContents()
$ws.Range("I12").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J12").Value = '01/23/20 08:24'
$ws.Range("K12").Value = '01/23/20 08:24'
$ws.Range("M12").Value = '$2,780 as of 1/23/2020 6:24:32 AM'
$ws.Range("N12").Value = 2800
$ws.Range("A13").Value = 'L678988'
$ws.Range("C13").Value = 'PAYELESS MARKET'
$ws.Range("E13").Value = 2860
$ws.Range("H13").Value = 45132.04188321759
$ws.Range("I13").ClearContents()
$ws.Range("J13").Value = '07/18/23 15:33'
$ws.Range("K13").Value = '07/18/23 15:33'
$ws.Range("M13").Value = '$2,880 as of 7/17/2023 5:44:43 PM'
$ws.Range("N13").Value = 2880
$ws.Range("H14").ClearContents()
$ws.Range("A15").Value = 'L474792'
$ws.Range("C15").Value = 'NICK SHELL SERVICE'
$ws.Range("E15").Value = 3580
$ws.Range("H15").Value = 45155.04188321759
$ws.Range("I15").ClearContents()
$ws.Range("J15").Value = '07/18/23 08:53'
$ws.Range("K15").Value = '07/18/23 08:53'
$ws.Range("M15").Value = '$3,580 as of 7/18/2023 6:53:29 AM'
$ws.Range("N15").Value = 3640
$ws.Range("A16").Value = 'L475182'
$ws.Range("C16").Value = 'LA ESQUINA DE ORO'
$ws.Range("E16").Value = 3800
$ws.Range("H16").ClearContents()
$ws.Range("I16").Value = 'ATM Inactive greater than 48 minutes'
$ws.Range("J16").Value = '09/16/20 16:57'
$ws.Range("K16").Value = '09/15/20 23:38'
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = '$3,800 as of 9/16/2020 1:28:00 PM'
$ws.Range("N16").Value = 3800
$ws.Range("A17").Value = 'L474817'
$ws.Range("C17").Value = 'SAFETY MARKET'
$ws.Range("E17").Value = 4580
$ws.Range("H17").Value = 45138.04188321759
$ws.Range("J17").Value = '07/18/23 15:08'
$ws.Range("K17").Value = '07/18/23 15:08'
$ws.Range("L17").Value = 120
$ws.Range("M17").Value = '$4,620 as of 7/18/2023 10:05:20 AM'
$ws.Range("N17").Value = 4620
$ws.Range("A18").Value = 'L476340'
$ws.Range("C18").Value = 'DONUT & SANDWICH'
$ws.Range("E18").Value = 4640
$ws.Range("H18").Value = 45145.04188321759
$ws.Range("J18").Value = '07/18/23 14:23'
$ws.Range("K18").Value = '07/18/23 14:23'
$ws.Range("M18").Value = '$4,700 as of 7/18/2023 8:03:40 AM'
$ws.Range("N18").Value = 4700
$ws.Range("A19").Value = 'L488595'
$ws.Range("C19").Value = 'N S MART'
$ws.Range("E19").Value = 5480
$ws.Range("H19").Value = 45285.04188321759
$ws.Range("I19").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J19").Value = '07/16/23 01:58'
$ws.Range("K19").Value = '07/16/23 01:58'
$ws.Range("L19").Value = 60
$ws.Range("M19").Value = '$5,480 as of 7/15/2023 11:58:38 PM'
$ws.Range("N19").Value = 5580
$ws.Range("A20").Value = 'LK236828'
$ws.Range("C20").Value = 'WORLDWIDE AUTOMOTIVE'
$ws.Range("E20").Value = 5480
$ws.Range("H20").Value = 45151.04188321759
$ws.Range("J20").Value = '07/17/23 20:02'
$ws.Range("K20").Value = '07/17/23 20:02'
$ws.Range("L20").Value = 80
$ws.Range("M20").Value = '$5,480 as of 7/17/2023 6:02:33 PM'
$ws.Range("N20").Value = 5500
$ws.Range("A21").Value = 'L688966'
$ws.Range("C21").Value = 'LACON MINI MART'
$ws.Range("E21").Value = 6400
$ws.Range("H21").ClearContents()
$ws.Range("I21").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J21").Value = '07/09/23 22:56'
$ws.Range("K21").Value = '07/09/23 15:28'
$ws.Range("L21").Value = 20
$ws.Range("M21").Value = '$6,400 as of 7/9/2023 1:28:46 PM'
$ws.Range("N21").Value = 6420
$ws.Range("A22").Value = 'L474761'
$ws.Range("C22").Value = 'BABS MARKET'
$ws.Range("E22").Value = 6740
$ws.Range("H22").Value = 45180.04188321759
$ws.Range("I22").ClearContents()
$ws.Range("J22").Value = '07/18/23 14:05'
$ws.Range("K22").Value = '07/18/23 14:05'
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = '$6,780 as of 7/16/2023 9:53:05 AM'
$ws.Range("N22").Value = 6780
$ws.Range("A23").Value = 'L475090'
$ws.Range("C23").Value = 'S.B. 2'
$ws.Range("E23").Value = 7900
$ws.Range("H23").Value = 45134.04188321759
$ws.Range("J23").Value = '07/18/23 12:20'
$ws.Range("K23").Value = '07/18/23 12:20'
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = '$7,900 as of 7/18/2023 10:20:54 AM'
$ws.Range("N23").Value = 8000
$ws.Range("A24").Value = 'LK923383'
$ws.Range("C24").Value = 'SAMYS PHONE CARDS'
$ws.Range("E24").Value = 8180
$ws.Range("H24").Value = 45133.04188321759
$ws.Range("J24").Value = '07/18/23 15:32'
$ws.Range("K24").Value = '07/18/23 15:32'
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = '$8,200 as of 7/17/2023 8:23:28 PM'
$ws.Range("N24").Value = 8200
$ws.Range("A25").Value = 'L697590'
$ws.Range("C25").Value = 'S B MARKET ST'
$ws.Range("E25").Value = 8780
$ws.Range("H25").ClearContents()
$ws.Range("I25").Value = 'ATM Inactive greater than 2000 minutes'
$ws.Range("J25").Value = '06/29/23 11:36'
$ws.Range("K25").Value = '06/29/23 11:36'
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = '$8,780 as of 6/29/2023 9:36:36 AM'
$ws.Range("N25").Value = 8800
$ws.Range("A26").Value = 'L697589'
$ws.Range("C26").Value = 'S B DISCOUNT MART'
$ws.Range("E26").Value = 11360
$ws.Range("H26").Value = 45131.04188321759
$ws.Range("I26").ClearContents()
$ws.Range("J26").Value = '07/18/23 15:06'
$ws.Range("K26").Value = '07/18/23 15:06'
$ws.Range("L26").Value = 40
$ws.Range("M26").Value = '$11,680 as of 7/18/2023 10:57:39 AM'
$ws.Range("N26").Value = 11540
$ws.Range("E27").Value = 16320
$ws.Range("H27").Value = 45129.04188321759
$ws.Range("J27").Value = '07/18/23 15:45'
$ws.Range("K27").Value = '07/18/23 15:45'
$ws.Range("M27").Value = '$16,600 as of 7/18/2023 11:33:29 AM'
$ws.Range("N27").Value = 16360
$ws.Range("E28").Value = 112300
